# Update "想去人数" (F column) and "最低票价" (G column) figures for the
# 合肥-漫展信息 workbook. The same underlying rows (5-13) are duplicated on
# both the "展览" sheet and the "全部类型" sheet, so both need updating.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F5").Value = 518
    $ws.Range("G5").Value = 80

    $ws.Range("F6").Value = 6921

    $ws.Range("F7").Value = 193

    $ws.Range("F8").Value = 155

    $ws.Range("F9").Value = 1044

    $ws.Range("F10").Value = 402

    $ws.Range("F11").Value = 138

    $ws.Range("F13").Value = 583
}
